$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Preconditions (E), Method Inputs (F), and Expected Result (G) columns
# for test case rows 7-13, as the rectangle test plan was completed with all test cases.

# Row 7
$ws.Range("E7").Value = 'None'
$ws.Range("F7").Value = 'color="red",                                                           length=5,                                                                width=6'
$ws.Range("G7").Value = 'Attributes set: _color="red",                   _length=5,                                                 -width=6'
$ws.Rows.Item(7).RowHeight = 72.75

# Row 8
$ws.Range("E8").Value = 'None'
$ws.Range("F8").Value = 'color=" ",                                                 length=5,                                                   width=6'
$ws.Range("G8").Value = 'ValueError: "Color cannot be blank."'

# Row 9
$ws.Range("E9").Value = 'None'
$ws.Range("F9").Value = 'color="red",                                       length="five",                                          width=6'
$ws.Range("G9").Value = 'ValueError:"Length must be numeric."'

# Row 10
$ws.Range("E10").Value = 'None'
$ws.Range("F10").Value = 'color="red",                                           length=5,                                           width="six"'
$ws.Range("G10").Value = 'ValueError: "Width must be numeric."'

# Row 11
$ws.Range("E11").Value = 'Valid instance'
$ws.Range("F11").Value = 'None'
$ws.Range("G11").Value = '"The shape color is red.
This rectangle has four sides with the lengths of 5, 6, 5 and 6 centimeters."'
$ws.Rows.Item(11).RowHeight = 91.5

# Row 12
$ws.Range("E12").Value = 'Valid instance'
$ws.Range("F12").Value = 'None'
$ws.Range("G12").Value = 'Correct area value(example, 12 for length=3, width=4) '

# Row 13
$ws.Range("E13").Value = 'Valid instance'
$ws.Range("F13").Value = 'None'
$ws.Range("G13").Value = 'Correct perimeter value(example, 14 for length=3, width=4)'
